# Apply the "Added many more features" edits to the Alkemor's Tower review.
#
# We replace whole-paragraph text (minus the trailing paragraph mark) via
# Range.Text assignment rather than Find/Replace, since several of the
# paragraphs share identical source text (e.g. the H1 title also appears
# verbatim near the end of the document) and Find/Replace's "smart quotes"
# autocorrect would otherwise mangle the straight apostrophes.

function Set-ParagraphText {
    param($Document, [int]$Index, [string]$OldText, [string]$NewText)

    $para = $Document.Paragraphs.Item($Index)
    $fullRange = $para.Range
    # Exclude the trailing paragraph-mark character from the range so we
    # only touch the visible text, preserving paragraph formatting.
    $textRange = $Document.Range($fullRange.Start, $fullRange.End - 1)

    if ($textRange.Text -ne $OldText) {
        throw "Paragraph $Index text mismatch. Expected '$OldText' but found '$($textRange.Text)'."
    }

    $textRange.Text = $NewText
}

$d = $word.ActiveDocument

# NOTE: this runtime's PowerShell does not bind named (-Param value)
# arguments correctly, so call Set-ParagraphText positionally:
#   Set-ParagraphText <Document> <Index> <OldText> <NewText>

Set-ParagraphText $d 1 `
    "Play Alkemor's Tower for Free - Review & Gameplay Mechanics" `
    "Play Alkemor's Tower for Free"

Set-ParagraphText $d 44 `
    "Exciting and unique gameplay mechanics" `
    "Exciting gameplay mechanics and special functions"

Set-ParagraphText $d 45 `
    "Special functions that provide advantageous situations and potentially excellent winnings" `
    "Unique theme with impressive graphics and animations"

Set-ParagraphText $d 46 `
    "Beautiful graphics and animations" `
    "Opportunities to win big with free spin bonuses"

Set-ParagraphText $d 47 `
    "Two free spin bonuses with higher winnings" `
    "Suitable for both casual players and high rollers"

Set-ParagraphText $d 49 `
    "RTP percentage is slightly lower than average" `
    "Medium level of volatility may not appeal to all players"

Set-ParagraphText $d 50 `
    "No jackpot feature" `
    "RTP percentage is slightly lower compared to some other slot games"

Set-ParagraphText $d 51 `
    "Play Alkemor's Tower for Free - Review & Gameplay Mechanics" `
    "Play Alkemor's Tower for Free"

Set-ParagraphText $d 52 `
    "Discover the magical world of Alkemor's Tower, a unique and exciting slot game. Learn how to play and trigger its special functions for better winnings. Play for free." `
    "Read our review of Alkemor's Tower and discover its exciting gameplay. Play for free and win big!"

Write-Output "Done."
